$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row that only contained the supervisor name value
# (B13:C13 = "4780627 - Ana Lucia Gabas Ferreira", no label in column A).
# Everything below shifts up by one row.
$ws.Rows(13).Delete()

# Update Objetivos: value (row 10, was long description text)
$ws.Range("B10:C10").Value = "4780627 - Ana Lucia Gabas Ferreira"

# Update Programa resumido: value (now row 13 after the deletion above)
$ws.Range("B13:C13").Value = "Semestral"

# Update Programa: value (now row 15 after the deletion above)
$ws.Range("B15:C15").Value = "01/01/2022"

# Update Método: value (now row 18 after the deletion above)
$ws.Range("B18:C18").Value = "4780627 - Ana Lucia Gabas Ferreira"

# Update Critério: value (now row 19 after the deletion above)
$ws.Range("B19:C19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."

# Update Norma de recuperação: value (now row 20 after the deletion above)
$ws.Range("B20:C20").Value = "A nota final será baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."

# Update Bibliografia: value (now row 21 after the deletion above)
$ws.Range("B21:C21").Value = "Devido às características da disciplina, não será oferecida recuperação."
